$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the player name
$ws.Name = "Shakib Al Hasan"

# A "matchNo" column is being inserted at the front of the table, pushing
# the original A:L columns to B:M. Use a real column insert so the sheet
# structure matches the diff (rather than just overwriting cells).
$ws.Columns.Item(1).Insert()

# Every column in this table stores numeric-looking values ("26", "104.00",
# "0", ...) as literal text (see the numberStoredAsText ignored-error on the
# source file), so force Text format across the whole used range before
# typing anything in, otherwise Excel will silently coerce them to numbers.
$ws.Range("A1:M7").NumberFormat = "@"

# Full header + data grid for the table, in final (post-insert) column order:
# matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr,
# opponentTeamName, venue, date, result
$grid = @(
    @("matchNo", "teamName", "batterName", "states", "runs", "balls", "fours", "sixes", "sr", "opponentTeamName", "venue", "date", "result"),
    @("10th", "Kolkata Knight Riders", "Shakib Al Hasan", "b Jamieson", "26", "25", "1", "1", "104.00", "Royal Challengers Bangalore", "Chennai", "April 18", "RCB won by 38 runs"),
    @("Qualifier", "Kolkata Knight Riders", "Shakib Al Hasan", "lbw b Ashwin", "0", "2", "0", "0", "0.00", "Delhi Capitals", "Sharjah", "October 13", "KKR won by 3 wickets (with 1 ball remaining)"),
    @("Eliminator", "Kolkata Knight Riders", "Shakib Al Hasan", "", "9", "6", "1", "0", "150.00", "Royal Challengers Bangalore", "Sharjah", "October 11", "KKR won by 4 wickets (with 2 balls remaining)"),
    @("Final", "Kolkata Knight Riders", "Shakib Al Hasan", "lbw b Jadeja", "0", "1", "0", "0", "0.00", "Chennai Super Kings", "Dubai (DSC)", "October 15", "Super Kings won by 27 runs"),
    @("3rd", "Kolkata Knight Riders", "Shakib Al Hasan", "c Abdul Samad b Kumar", "3", "5", "0", "0", "60.00", "Sunrisers Hyderabad", "Chennai", "April 11", "KKR won by 10 runs"),
    @("5th", "Kolkata Knight Riders", "Shakib Al Hasan", "c Yadav b KH Pandya", "9", "9", "1", "0", "100.00", "Mumbai Indians", "Chennai", "April 13", "Mumbai won by 10 runs")
)

$r = 1
foreach ($row in $grid) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
